$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three e-mail addresses in column A (values only - underlying
# hyperlink targets are left untouched, matching the source edit).
$ws.Range("A1").Value = "darsh3@gmail.com"
$ws.Range("A2").Value = "sanj3@gmail.com"
$ws.Range("A3").Value = "harshi3@gmail.com"

# Move/record the active selection on the sheet to D2, as in the saved file.
$ws.Range("D2").Select() | Out-Null
